$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.869.70"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "3.113.65"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.36"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.55"
$ws.Range("E6").Value = "  +1.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.110.02"
$ws.Range("E8").Value = "  +0.35%  "

$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.44"
$ws.Range("E10").Value = "  -3.86%  "

$ws.Range("E11").Value = "  -1.96%  "

$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("E13").Value = "  -2.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.39"
$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").Value = "3.628.66"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").Value = "66.814.41"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("D19").Value = "3.111.16"
$ws.Range("E19").Value = "  +0.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.52"
$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "477.80"
$ws.Range("E21").Value = "  +1.34%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.01"
$ws.Range("E22").Value = "  +6.09%  "

$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.715"
$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.59"
$ws.Range("E24").Value = "  +5.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.98"
$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("E26").Value = "  -0.96%  "

$ws.Range("E27").Value = "  -1.87%  "

$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.42"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.93"
$ws.Range("E30").Value = "  -2.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.65"
$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("E33").Value = "  -2.12%  "

$ws.Range("D34").Value = "0.0₃0944"
$ws.Range("E34").Value = "  -7.54%  "

$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.976"
$ws.Range("E37").Value = "  -3.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.09"
$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("E40").Value = "  -0.52%  "

$ws.Range("E41").Value = "  -2.24%  "

$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.68"
$ws.Range("E43").Value = "  -0.63%  "

$ws.Range("D44").Value = "2.808.86"
$ws.Range("E44").Value = "  +1.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0356"
$ws.Range("E45").Value = "  -2.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "379.75"
$ws.Range("E46").Value = "  -4.34%  "

$ws.Range("E47").Value = "  -11.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.31"
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.99"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("E51").Value = "  -1.98%  "
